$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test "username" credential for row 2 (was Admin@gmail.com)
$ws.Range("A2").Value = "admin_ankush@gmail.com"

# Row 2's password cell becomes "Test@admin" and - because it looks like an
# email/short token - picked up a hyperlink (mirroring the Excel auto-link
# behaviour seen in the authored workbook). Set the text first, then attach
# the hyperlink so no extra "display text" override is written.
$ws.Range("B2").Value = "Test@admin"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:test@admin.com")

# Update the test "username" credential for row 3 (was Admin@gmail.com)
$ws.Range("A3").Value = "admin_new@gmail.com"

# Row 3's password ("Admin@1234") is left as-is.
